$d = $word.ActiveDocument

$replacements = @(
    @("766÷8=", "892÷5="),
    @("916÷6=", "956÷3="),
    @("375÷4=", "644÷2="),
    @("825÷5=", "545÷6="),
    @("936÷6=", "204÷4="),
    @("279÷7=", "250÷5="),
    @("851÷6=", "788÷8="),
    @("351÷4=", "781÷7="),
    @("152÷3=", "428÷3="),
    @("692÷6=", "672÷4="),
    @("432÷5=", "514÷9="),
    @("772÷2=", "837÷7="),
    @("791÷2=", "114÷8="),
    @("498÷4=", "457÷2="),
    @("786÷6=", "865÷4="),
    @("116÷9=", "515÷6="),
    @("583÷9=", "398÷3="),
    @("512÷2=", "802÷4="),
    @("424÷2=", "671÷8="),
    @("614÷4=", "150÷5="),
    @("538÷2=", "471÷2="),
    @("673÷9=", "450÷4="),
    @("756÷2=", "880÷5="),
    @("653÷6=", "634÷9="),
    @("289÷8=", "330÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
